$d = $word.ActiveDocument

# 1) "В подвижном блоке предлагается..." -> "В блоке движения предлагается..."
#    ("подвижном " removed, "движения " inserted after "блоке ")
$d.Content.Find.Execute(
    "В подвижном блоке предлагается использовать д", $true, $false, $false, $false, $false,
    $true, 1, $false, "В блоке движения предлагается использовать д", 2)

# 2) "...является допустимым (подключим через резистор 1 кОм), а максимальный
#    допустимый ток сток-исток составляет 0.825мА, что превышает в 1.65 раза
#    максимальный возможный ток в цепи коллекторного двигателя."
#    -> add a clarifying note after "допустимым", turn the leading "(" before
#    "подключим" into a comma, and append a follow-up question at the end.
$d.Content.Find.Execute(
    "мым (подключим через резистор 1 кОм), а максимальный допустимый ток сток-исток составляет 0.825мА, что превышает в 1.65 раза максимальный возможный ток в цепи коллекторного двигателя.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "мым (максимальное не указано в даташите!), подключим через резистор 1 кОм, а максимальный допустимый ток сток-исток составляет 0.825мА, что превышает в 1.65 раза максимальный возможный ток в цепи коллекторного двигателя. (Если этот вариант неверно подобран, то как подключить транзистор с меньшим 3.3В упарвляющим напряжением? Допустимо ли подключить через делитель?)",
    2)
